# Insert a new price-report row for "Pepino ensalada" (Terminal Hortofrutícola
# Agro Chillán) dated 2022-01-28, pushing the existing row 98 (and everything
# below it) down by one row. This grows the used range from A1:R198 to
# A1:R199, matching the weekly refresh of the consolidated sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 98..198 down to 99..199 and leave a blank row 98 (formatting of
# the surrounding rows - e.g. the date-style on column D - is carried over
# automatically by Excel's row insert).
$ws.Rows.Item(98).Insert()

# Populate the newly inserted row 98 with the new observation.
$ws.Cells.Item(98, 1).Value  = 7
$ws.Cells.Item(98, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(98, 3).Value  = "Ñuble"
$ws.Cells.Item(98, 4).Value  = Get-Date -Year 2022 -Month 1 -Day 28 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(98, 5).Value  = 16
$ws.Cells.Item(98, 6).Value  = 100112043
$ws.Cells.Item(98, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(98, 8).Value  = "Sin especificar"
$ws.Cells.Item(98, 9).Value  = "Primera"
$ws.Cells.Item(98, 10).Value = 100
$ws.Cells.Item(98, 11).Value = 8000
$ws.Cells.Item(98, 12).Value = 8500
$ws.Cells.Item(98, 13).Value = 8250
$ws.Cells.Item(98, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(98, 15).Value = "Región del Maule"
$ws.Cells.Item(98, 16).Value = 103
$ws.Cells.Item(98, 17).Value = 80
$ws.Cells.Item(98, 18).Value = "Hortaliza"
